$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.029899180750632
$ws.Range("D2").Value = 1.035170097181239
$ws.Range("E2").Value = 1.039133933876635
$ws.Range("F2").Value = 1.04881679467445
$ws.Range("I2").Value = 1.037941592842828
$ws.Range("J2").Value = 1.035043404736897
$ws.Range("K2").Value = 1.037967470451768
$ws.Range("L2").Value = 1.041919986588995
$ws.Range("M2").Value = 1.051575583816734
$ws.Range("N2").Value = 1.015685007260422

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.03071605375388
$ws.Range("D3").Value = 1.035642765531242
$ws.Range("E3").Value = 1.039870692091212
$ws.Range("F3").Value = 1.049661549215666
$ws.Range("I3").Value = 1.038105440698532
$ws.Range("J3").Value = 1.035502182042361
$ws.Range("K3").Value = 1.038250374178817
$ws.Range("L3").Value = 1.042467094668118
$ws.Range("M3").Value = 1.052232375384859
$ws.Range("N3").Value = 1.015837681183171

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.031245202407314
$ws.Range("D4").Value = 1.03594872305099
$ws.Range("E4").Value = 1.040348337013234
$ws.Range("F4").Value = 1.050209191501308
$ws.Range("I4").Value = 1.038210255010776
$ws.Range("J4").Value = 1.035798955398388
$ws.Range("K4").Value = 1.038432817429971
$ws.Range("L4").Value = 1.042821354721283
$ws.Range("M4").Value = 1.052657745837523
$ws.Range("N4").Value = 1.015936415760608

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.031467792995269
$ws.Range("D5").Value = 1.036077371878461
$ws.Range("E5").Value = 1.04054935565645
$ws.Range("F5").Value = 1.050439664799052
$ws.Range("I5").Value = 1.038254029283011
$ws.Range("J5").Value = 1.035923696993636
$ws.Range("K5").Value = 1.038509368190206
$ws.Range("L5").Value = 1.042970342723245
$ws.Range("M5").Value = 1.052836661415446
$ws.Range("N5").Value = 1.015977909991341

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.031505174889467
$ws.Range("D6").Value = 1.036098973950678
$ws.Range("E6").Value = 1.040583120224976
$ws.Range("F6").Value = 1.050478376557894
$ws.Range("I6").Value = 1.038261362163791
$ws.Range("J6").Value = 1.035944640333582
$ws.Range("K6").Value = 1.038522212660613
$ws.Range("L6").Value = 1.042995361768583
$ws.Range("M6").Value = 1.052866707346886
$ws.Range("N6").Value = 1.015984876233466

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.031248176139428
$ws.Range("D7").Value = 1.035950441969927
$ws.Range("E7").Value = 1.040351022184991
$ws.Range("F7").Value = 1.050212270138675
$ws.Range("I7").Value = 1.03821084106405
$ws.Range("J7").Value = 1.035800622289438
$ws.Range("K7").Value = 1.038433840889656
$ws.Range("L7").Value = 1.042823345283834
$ws.Range("M7").Value = 1.052660136165365
$ws.Range("N7").Value = 1.015936970263297

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.030175126278562
$ws.Range("D8").Value = 1.035329813836348
$ws.Range("E8").Value = 1.039382734574588
$ws.Range("F8").Value = 1.049102069030974
$ws.Range("I8").Value = 1.03799721520714
$ws.Range("J8").Value = 1.035198467763302
$ws.Range("K8").Value = 1.038063205521666
$ws.Range("L8").Value = 1.042104832809476
$ws.Range("M8").Value = 1.051797469400993
$ws.Range("N8").Value = 1.015736615375359

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.028288770659261
$ws.Range("D9").Value = 1.034237111286104
$ws.Range("E9").Value = 1.03768356427295
$ws.Range("F9").Value = 1.047153727281261
$ws.Range("I9").Value = 1.037611580650261
$ws.Range("J9").Value = 1.034136790546875
$ws.Range("K9").Value = 1.03740545057805
$ws.Range("L9").Value = 1.0408406598107
$ws.Range("M9").Value = 1.050280344366199
$ws.Range("N9").Value = 1.015383158522442

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.027034323168298
$ws.Range("D10").Value = 1.033509378643423
$ws.Range("E10").Value = 1.036555653709336
$ws.Range("F10").Value = 1.045860314231421
$ws.Range("I10").Value = 1.037348359780937
$ws.Range("J10").Value = 1.033428675749359
$ws.Range("K10").Value = 1.036963899579962
$ws.Range("L10").Value = 1.039999268172545
$ws.Range("M10").Value = 1.049271049071549
$ws.Range("N10").Value = 1.015147275114556

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.026491895408498
$ws.Range("D11").Value = 1.033194461197344
$ws.Range("E11").Value = 1.036068434671992
$ws.Range("F11").Value = 1.045301577547074
$ws.Range("I11").Value = 1.037232939254775
$ws.Range("J11").Value = 1.033121991998951
$ws.Range("K11").Value = 1.036771996607531
$ws.Range("L11").Value = 1.039635284547152
$ws.Range("M11").Value = 1.048834537689519
$ws.Range("N11").Value = 1.015045082547776

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.026290528905962
$ws.Range("D12").Value = 1.033077518199142
$ws.Range("E12").Value = 1.035887638082301
$ws.Range("F12").Value = 1.045094238187076
$ws.Range("I12").Value = 1.037189850855352
$ws.Range("J12").Value = 1.033008067466675
$ws.Range("K12").Value = 1.036700610175948
$ws.Range("L12").Value = 1.039500138070927
$ws.Range("M12").Value = 1.048672477804274
$ws.Range("N12").Value = 1.015007116106457

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.026333717458626
$ws.Range("D13").Value = 1.033102601418318
$ws.Range("E13").Value = 1.035926411467305
$ws.Range("F13").Value = 1.045138704085684
$ws.Range("I13").Value = 1.037199103226935
$ws.Range("J13").Value = 1.033032505030994
$ws.Range("K13").Value = 1.036715927543783
$ws.Range("L13").Value = 1.039529125015315
$ws.Range("M13").Value = 1.048707236557346
$ws.Range("N13").Value = 1.015015260372184

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.02647524802829
$ws.Range("D14").Value = 1.033184794001429
$ws.Range("E14").Value = 1.036053486322049
$ws.Range("F14").Value = 1.045284434712518
$ws.Range("I14").Value = 1.03722938195699
$ws.Range("J14").Value = 1.033112575128179
$ws.Range("K14").Value = 1.036766097919702
$ws.Range("L14").Value = 1.039624112202397
$ws.Range("M14").Value = 1.048821140129354
$ws.Range("N14").Value = 1.015041944381669

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.026562465015027
$ws.Range("D15").Value = 1.033235439779495
$ws.Range("E15").Value = 1.036131805011917
$ws.Range("F15").Value = 1.04537425076744
$ws.Range("I15").Value = 1.037248009076249
$ws.Range("J15").Value = 1.03316190786563
$ws.Range("K15").Value = 1.036796995653639
$ws.Range("L15").Value = 1.039682644043848
$ws.Range("M15").Value = 1.048891330515731
$ws.Range("N15").Value = 1.015058384291148

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.027070338391006
$ws.Range("D16").Value = 1.03353028296703
$ws.Range("E16").Value = 1.036588013724436
$ws.Range("F16").Value = 1.045897423724792
$ws.Range("I16").Value = 1.037355989510533
$ws.Range("J16").Value = 1.033449028072594
$ws.Range("K16").Value = 1.036976620733987
$ws.Range("L16").Value = 1.040023431952615
$ws.Range("M16").Value = 1.049300030001033
$ws.Range("N16").Value = 1.015154056211292

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.027389117616848
$ws.Range("D17").Value = 1.033715284072692
$ws.Range("E17").Value = 1.036874496981962
$ws.Range("F17").Value = 1.046225951069465
$ws.Range("I17").Value = 1.037423336772401
$ws.Range("J17").Value = 1.033629114401446
$ws.Range("K17").Value = 1.037089106053258
$ws.Range("L17").Value = 1.040237292490962
$ws.Range("M17").Value = 1.049556536870269
$ws.Range("N17").Value = 1.01521405468566

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.02757512900234
$ws.Range("D18").Value = 1.033823210827105
$ws.Range("E18").Value = 1.037041710933479
$ws.Range("F18").Value = 1.046417702705873
$ws.Range("I18").Value = 1.037462479891318
$ws.Range("J18").Value = 1.03373414931658
$ws.Range("K18").Value = 1.037154648391049
$ws.Range("L18").Value = 1.040362066864725
$ws.Range("M18").Value = 1.049706203057275
$ws.Range("N18").Value = 1.01524904558806

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.027638566383956
$ws.Range("D19").Value = 1.033860014179533
$ws.Range("E19").Value = 1.037098745697386
$ws.Range("F19").Value = 1.046483106570525
$ws.Range("I19").Value = 1.03747580301164
$ws.Range("J19").Value = 1.033769962384594
$ws.Range("K19").Value = 1.037176984984729
$ws.Range("L19").Value = 1.04040461727747
$ws.Range("M19").Value = 1.049757243799329
$ws.Range("N19").Value = 1.015260975686618

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.027354908074692
$ws.Range("D20").Value = 1.033695433252303
$ws.Range("E20").Value = 1.036843748325867
$ws.Range("F20").Value = 1.046190690023403
$ws.Range("I20").Value = 1.037416125460887
$ws.Range("J20").Value = 1.033609793494705
$ws.Range("K20").Value = 1.037077044513216
$ws.Range("L20").Value = 1.040214343840574
$ws.Range("M20").Value = 1.049529010919287
$ws.Range("N20").Value = 1.015207617950918

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.026433567606078
$ws.Range("D21").Value = 1.033160589452685
$ws.Range("E21").Value = 1.036016061008472
$ws.Range("F21").Value = 1.045241515128805
$ws.Range("I21").Value = 1.037220471583296
$ws.Range("J21").Value = 1.033088996705439
$ws.Range("K21").Value = 1.036751326885811
$ws.Range("L21").Value = 1.039596139361409
$ws.Range("M21").Value = 1.048787596147446
$ws.Range("N21").Value = 1.015034086811247

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.025854952327236
$ws.Range("D22").Value = 1.03282449435757
$ws.Range("E22").Value = 1.035496693414741
$ws.Range("F22").Value = 1.04464589121803
$ws.Range("I22").Value = 1.037096206511183
$ws.Range("J22").Value = 1.032761502365936
$ws.Range("K22").Value = 1.036545927724928
$ws.Range("L22").Value = 1.039207758485817
$ws.Range("M22").Value = 1.048321902258295
$ws.Range("N22").Value = 1.014924937317436

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.026161623344719
$ws.Range("D23").Value = 1.033002646776132
$ws.Range("E23").Value = 1.035771921514165
$ws.Range("F23").Value = 1.044961532081484
$ws.Range("I23").Value = 1.037162199938217
$ws.Range("J23").Value = 1.032935117500372
$ws.Range("K23").Value = 1.036654870857275
$ws.Range("L23").Value = 1.03941361676676
$ws.Range("M23").Value = 1.048568730934474
$ws.Range("N23").Value = 1.014982803504924

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.027370365669351
$ws.Range("D24").Value = 1.033704402925931
$ws.Range("E24").Value = 1.036857641970741
$ws.Range("F24").Value = 1.046206622579577
$ws.Range("I24").Value = 1.037419384372984
$ws.Range("J24").Value = 1.033618523800719
$ws.Range("K24").Value = 1.037082494815636
$ws.Range("L24").Value = 1.040224713245663
$ws.Range("M24").Value = 1.049541448557206
$ws.Range("N24").Value = 1.015210526450385

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.028775895732426
$ws.Range("D25").Value = 1.034519479788303
$ws.Range("E25").Value = 1.038121990410661
$ws.Range("F25").Value = 1.047656462955322
$ws.Range("I25").Value = 1.037712360993578
$ws.Range("J25").Value = 1.034411323287664
$ws.Range("K25").Value = 1.037576038828366
$ws.Range("L25").Value = 1.041167240174548
$ws.Range("M25").Value = 1.050672191081885
$ws.Range("N25").Value = 1.01547458086921
